$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell K1, copying style from J1 (bold/centered header style)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "intervention_type"

# Fill in intervention_type values for each row
$ws.Range("K2").Value = "OTHER"
$ws.Range("K3").Value = "DRUG"
$ws.Range("K4").Value = "DEVICE"
$ws.Range("K5").Value = "BEHAVIORAL"
$ws.Range("K6").Value = "DEVICE"
